# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the freshly generated site snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" sheet (column F)
$exhibitUpdates = @{
    2  = 274
    4  = 13684
    5  = 1344
    6  = 259
    7  = 44
    8  = 99
    10 = 252
    11 = 485
    12 = 10
    13 = 81
    18 = 5670
    19 = 118
    20 = 76
    21 = 970
    22 = 55
    23 = 48
    24 = 138
    25 = 189
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new value for the "全部类型" sheet (column F)
# Note: F4 differs slightly from the "展览" sheet (13685 vs 13684) in the
# source data, matching the upstream snapshot exactly.
$allUpdates = @{
    2  = 274
    4  = 13685
    5  = 1344
    6  = 259
    7  = 44
    8  = 99
    10 = 252
    11 = 485
    12 = 10
    13 = 81
    18 = 5670
    19 = 118
    20 = 76
    21 = 970
    22 = 55
    23 = 48
    24 = 138
    25 = 189
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
